# Updated cryptos list on Tue Aug  1 15:24:57 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.814.93"
$ws.Range("E2").Value = "  -1.75%  "

$ws.Range("D3").Value = "1.829.23"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6888"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07640"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3035"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07781"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "92.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").Value = "1.832.76"
$ws.Range("E13").Value = "  -1.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.073"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6756"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.440"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008215"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.49%  "

$ws.Range("D18").Value = "28.832.53"
$ws.Range("E18").Value = "  -1.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.98%  "

$ws.Range("D20").Value = "2.073.24"
$ws.Range("E20").Value = "  -2.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.408"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.02%  "

$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1481"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.708"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.538"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.209"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.146"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.182"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05077"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7692"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.839"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.133"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.695"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01854"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("D39").Value = "1.240.73"
$ws.Range("E39").Value = "  -2.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.701"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9539"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.83"
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.671"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.03%  "

$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("D47").Value = "1.974.49"
$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("E48").Value = "  -0.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.735"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.904"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.38%  "

